# Implement the new match logic (partially)
#
# The old "match" row (A42:F42) is kept but moved down to row 43, and a new
# "match_decide" CODE row is written in its place at row 42. Three more rows
# (44-46: "no_match" + two alternative-path options) are appended right after
# the (now relocated) "match" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: insert 4 new blank rows starting at row 43 -------------
$ws.Rows("43:46").Insert()

# --- 2. Row 43 = the old "match" row content, with two cells updated ------
#     (D: new recommendation text, F: reference to the new [matched] var
#      instead of the raw MATCH(...) formula string)
$ws.Range("A43").Value = "match"
$ws.Range("B43").Value = "Match"
$ws.Range("C43").Value = "PN"
$ws.Range("D43").Value = "Folgende Arzneiformen passen optimal zu deinen Bedürfnissen:"
$ws.Range("E43").Value = "JUMP(Night)"
$ws.Range("F43").Value = "[matched]"

# --- 3. Row 42 becomes the new "match_decide" CODE row --------------------
$ws.Range("A42").Value = "match_decide"
$ws.Range("B42").ClearContents()
$ws.Range("C42").Value = "CODE"
$ws.Range("D42").ClearContents()
$ws.Range("E42").Value = "matched = MATCH([api_solver], [water], [swallow], [transport], [fly], [single]), SAVE(matched), IF (ROWS([matched]) == 0) THEN GO(no_match) ELSE GO(match)"
$ws.Range("F42").ClearContents()

# --- 4. Row 44 = new "no_match" question row -------------------------------
$ws.Range("A44").Value = "no_match"
$ws.Range("B44").Value = "No Match"
$ws.Range("C44").Value = "Q"
$ws.Range("D44").Value = "Leider gibt es mit dem Wirkstoff/Extrakt {SELECT Name FROM API WHERE ID == [api_solver]},den du ausgewählt hast (show chosen API), keine passenden Medikamente."

# --- 5. Row 45 = first alternative answer (alternative API) ---------------
$ws.Range("C45").Value = "A"
$ws.Range("D45").Value = "Zeig mir Produkte in der passenden Arzneiformen aber alternativen Wirkstoffen"
$ws.Range("E45").Value = "FINISH()"

# --- 6. Row 46 = second alternative answer (alternative Arzneiform) -------
$ws.Range("C46").Value = "A"
$ws.Range("D46").Value = "Zeig mir Produkte mit dem ausgewählten Wirkstoff aber in einer anderen Arzneiform"
$ws.Range("E46").Value = "FINISH()"

# --- 7. Cosmetic touch-ups carried in the diff (column widths / row -------
#     heights / dimension) -- best-effort visual parity, not required for
#     the underlying content change.
$ws.Columns(1).ColumnWidth = 14.7
$ws.Columns(5).ColumnWidth = 34.12
$ws.Rows(42).RowHeight = 68.65
$ws.Rows(43).RowHeight = 14.2
$ws.Rows(44).RowHeight = 41.9
$ws.Rows(45).RowHeight = 28.35
$ws.Rows(46).RowHeight = 28.35

$ws.Range("E39").Select()
